$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91: new date entry (09/11/2021) with DONE note
$ws.Range("A91").Value = "09/11/2021"
$ws.Range("B91").Value = "Explored on distinct layers of android multimedia framework"

# Row 92: PROGRESS / TO-DO continuation
$ws.Range("B92").Value = "studied more on stagefright"
$ws.Range("C92").Value = "Updating the notes"
$ws.Range("D92").Value = "Listed the doubts , have to elaborate "

# Row 93
$ws.Range("B93").Value = "Listed the links used for self study"
$ws.Range("D93").Value = "Revision of C-DS-OS concepts"

# Row 94
$ws.Range("B94").Value = "Discussed and guided the new teammate on V4L2 testapps"

# Update selection / view to match final state
$ws.Range("D93").Select()
